$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Grade tables now also in stripped form for easier input: mark every
# data row (10-50) in the leading "stripped" column (A) with "X", just
# like the example already shown in row 9.
for ($r = 10; $r -le 50; $r++) {
    $ws.Cells.Item($r, 1).Value = "X"
}

# Bring the view back to the top of the frozen pane instead of leaving it
# scrolled out to the far right/bottom where it had been left selected.
$ws.Range("A1").Select()
